$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = "0.741 (0.726 ± 0.016)"
    "C2" = "00:04:54 (00:10:22 ± 00:06:53)"
    "D2" = "00:00:03 (00:00:09 ± 00:00:04)"

    "B3" = "0.799 (0.742 ± 0.026)"
    "C3" = "00:01:28 (00:03:38 ± 00:01:36)"
    "D3" = "00:00:00 (00:00:00 ± 00:00:00)"

    "B4" = "0.716 (0.679 ± 0.024)"
    "C4" = "00:01:16 (00:01:40 ± 00:00:21)"
    "D4" = "00:00:01 (00:00:01 ± 00:00:00)"

    "B5" = "0.800 (0.625 ± 0.151)"
    "C5" = "00:05:07 (00:05:13 ± 00:00:04)"
    "D5" = "00:00:01 (00:00:02 ± 00:00:02)"

    "B6" = "0.796 (0.749 ± 0.021)"
    "C6" = "00:04:56 (00:05:00 ± 00:00:02)"
    "D6" = "00:00:01 (00:00:05 ± 00:00:01)"

    "B9" = "0.772 (0.723 ± 0.024)"
    "C9" = "00:05:01 (00:05:04 ± 00:00:02)"
    "D9" = "00:00:00 (00:00:00 ± 00:00:00)"

    "B11" = "0.607 (0.472 ± 0.143)"
    "C11" = "00:05:05 (00:05:06 ± 00:00:00)"
    "D11" = "00:00:00 (00:00:00 ± 00:00:00)"

    "B12" = "0.773 (0.773 ± 0.000)"
    "C12" = "00:02:04 (00:02:04 ± 00:00:00)"
    "D12" = "00:00:00 (00:00:00 ± 00:00:00)"

    "B13" = "0.299 (0.235 ± 0.040)"
    "C13" = "00:00:19 (00:00:19 ± 00:00:00)"
    "D13" = "00:00:01 (00:00:01 ± 00:00:00)"

    "B14" = "0.749 (0.692 ± 0.022)"
    "C14" = "00:02:08 (00:02:25 ± 00:00:10)"
    "D14" = "00:00:00 (00:00:00 ± 00:00:00)"

    "B15" = "0.781 (0.726 ± 0.030)"
    "C15" = "00:00:51 (00:04:15 ± 00:01:13)"
    "D15" = "00:00:00 (00:00:00 ± 00:00:00)"

    "B16" = "0.785 (0.712 ± 0.026)"
    "C16" = "00:09:31 (00:11:36 ± 00:01:32)"
    "D16" = "00:00:00 (00:00:00 ± 00:00:00)"

    "B17" = "0.753 (0.704 ± 0.024)"
    "C17" = "00:05:03 (00:05:56 ± 00:00:35)"
    "D17" = "00:00:00 (00:00:00 ± 00:00:00)"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
